# Add a new "Time zone support" bullet (indent level 2 / lvl="1") right
# before the existing "SQL 2016 Service Pack 1" bullet on the slide that
# lists the free SQL Server Express features.

$p = $ppt.ActivePresentation

# Locate the slide/shape that contains the "SQL 2016 Service Pack 1" bullet.
$targetSlide = $null
$targetShape = $null
foreach ($sl in $p.Slides) {
    foreach ($sh in $sl.Shapes) {
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            if ($sh.TextFrame.TextRange.Text -like "*SQL 2016 Service Pack 1*") {
                $targetSlide = $sl
                $targetShape = $sh
            }
        }
    }
}

$tr = $targetShape.TextFrame.TextRange
$paraCount = $tr.Paragraphs().Count

# Find the paragraph index of "SQL 2016 Service Pack 1" - the new bullet is
# inserted immediately before it.
$insertBeforeIndex = -1
for ($i = 1; $i -le $paraCount; $i++) {
    $paraText = $tr.Paragraphs($i, 1).Text.TrimEnd("`r")
    if ($paraText -eq "SQL 2016 Service Pack 1") {
        $insertBeforeIndex = $i
        break
    }
}

$targetPara = $tr.Paragraphs($insertBeforeIndex, 1)

# Type the new bullet as three runs ("Time " / "zone suppor" / "t"),
# inserting back-to-front so each InsertBefore call lands right before the
# target paragraph and the pieces end up in reading order.
$rT = $targetPara.InsertBefore("t`r")
$rZoneSupport = $rT.InsertBefore("zone suppor")
$rTime = $rZoneSupport.InsertBefore("Time ")

# Re-fetch the freshly created paragraph and indent it one level, matching
# its sibling bullets ("Stretch DB", "Temporal tables", ...).
$newPara = $tr.Paragraphs($insertBeforeIndex, 1)
$newPara.IndentLevel = 2
